$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (H) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 404
$wsOff.Range("C2").Value = 287
$wsOff.Range("D2").Value = 112
$wsOff.Range("E2").Value = 54
$wsOff.Range("G2").Value = 3

# --- DEF sheet: row 2 (H) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 472
$wsDef.Range("C2").Value = 315
$wsDef.Range("D2").Value = 104
$wsDef.Range("E2").Value = 37
$wsDef.Range("G2").Value = 10
